$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update numeric values for the RandomForest imputation results (column A-D)
# per the commit "Update Name of Algo" data refresh.

$ws.Range("D2").Value = -7.413799999999993
$ws.Range("D5").Value = -8.184700000000007
$ws.Range("D6").Value = -7.954799999999999
$ws.Range("D8").Value = -8.735400000000002
$ws.Range("B11").Value = 6.1002
$ws.Range("A12").Value = -21.40519999999999
$ws.Range("C14").Value = -13.54739999999999
$ws.Range("D17").Value = -8.113799999999992
$ws.Range("C19").Value = -12.89330000000001
$ws.Range("B23").Value = 8.7441
$ws.Range("C24").Value = -13.34579999999998
$ws.Range("A27").Value = -21.8891
$ws.Range("D27").Value = -7.777299999999999
$ws.Range("B28").Value = 5.884800000000002
$ws.Range("A32").Value = -21.02609999999999
$ws.Range("B32").Value = 6.596900000000002
$ws.Range("B34").Value = 9.449000000000003
$ws.Range("A36").Value = -19.832
$ws.Range("A38").Value = -19.62749999999999
$ws.Range("C38").Value = -10.56250000000001
$ws.Range("C41").Value = -12.38810000000001
$ws.Range("B42").Value = 10.66089999999999
$ws.Range("A46").Value = -21.9753
$ws.Range("B49").Value = 5.3459
$ws.Range("C52").Value = -11.1626
$ws.Range("A54").Value = -21.84440000000001
$ws.Range("B54").Value = 4.9963
$ws.Range("A55").Value = -21.50829999999999
$ws.Range("D55").Value = -7.860600000000002
$ws.Range("A56").Value = -22.04560000000001
$ws.Range("A67").Value = -21.54799999999997
$ws.Range("A69").Value = -21.60469999999997
$ws.Range("D70").Value = -6.518299999999999
$ws.Range("A72").Value = -22.1323
$ws.Range("C72").Value = -11.9311
$ws.Range("B78").Value = 9.627800000000001
$ws.Range("C78").Value = -12.3019
$ws.Range("B80").Value = 9.476499999999998
$ws.Range("D80").Value = -7.294899999999993
$ws.Range("A83").Value = -21.61229999999999
$ws.Range("C83").Value = -13.06239999999998
$ws.Range("C85").Value = -14.0207
$ws.Range("A86").Value = -21.73509999999999
$ws.Range("C86").Value = -12.8637
$ws.Range("C90").Value = -9.968500000000006
$ws.Range("A91").Value = -20.74879999999999
$ws.Range("A93").Value = -21.30100000000002
$ws.Range("D95").Value = -7.491400000000006
$ws.Range("C96").Value = -10.1239
$ws.Range("B97").Value = 6.393799999999996
$ws.Range("D98").Value = -8.252900000000006
$ws.Range("A99").Value = -21.7621
$ws.Range("B99").Value = 5.725199999999996
$ws.Range("B101").Value = 4.956999999999997
$ws.Range("D102").Value = -7.965
$ws.Range("C103").Value = -13.81019999999999
$ws.Range("A104").Value = -21.6828
